$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")

# Rows 107-118 (probateFormsRW06.feature) flip SmokeTest from "Yes" to "No"
for ($r = 107; $r -le 118; $r++) {
    $ws.Cells.Item($r, 4).Value = "No"
}

# Rows 119-124 (probateFormsRWxx.feature) flip SmokeTest from "No" to "Yes"
for ($r = 119; $r -le 124; $r++) {
    $ws.Cells.Item($r, 4).Value = "Yes"
}

# Update the visible window / selection to match the edited area
$ws.Range("C124").Select()
$ws.Application.ActiveWindow.ScrollRow = 107
